# Sablefish_Input.xlsx: "Updated to include length bins"
#
# Adds a new "Length_Bins" worksheet (positioned right after "Sim" and
# before "Growth") that holds the model's length-bin labels, derived from
# the General sheet's parameters:
#   n.length = 30, first.len = 41, len.incr = 2  ->  41, 43, 45, ..., 99
#
# The workbook's active sheet/selection (Maturity!D9) is restored afterward
# so inserting the new sheet doesn't change which tab is focused.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately before "Growth" (i.e. right after "Sim").
$growthSheet = $wb.Worksheets.Item("Growth")
$lenBins = $wb.Worksheets.Add($growthSheet)
$lenBins.Name = "Length_Bins"

# Fill column A with the length-bin values: first.len=41, len.incr=2, n.length=30.
$firstLen = 41
$lenIncr = 2
$nLength = 30

$value = $firstLen
for ($row = 1; $row -le $nLength; $row++) {
    $lenBins.Cells.Item($row, 1).Value = $value
    $value = $value + $lenIncr
}

$lenBins.Range("D9").Select()

# Keep "Maturity" as the active/selected sheet (as in the authored workbook).
$maturitySheet = $wb.Worksheets.Item("Maturity")
$maturitySheet.Activate()
$maturitySheet.Range("D9").Select()
